$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.511.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4526"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3803"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07827"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.146"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9991"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.407"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.837.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +16.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001091"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06405"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.399"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5425"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "28.568.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.296"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.380"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.043.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.219"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.919"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09332"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.666"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02364"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6714"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06321"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.244"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.202"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.410"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6189"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.780"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.068"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("E51").Value = "  -0.61%  "
